# Adicionando Item do Fator de Ajuste - Taxas de Frequência e Gravidade
#
# Adds a new row (33) to the "Planilha1" worksheet of the Alteracoes_Casos
# workbook, describing a new "Caso" about the "Fator de Ajuste" for the
# Taxas de Frequência e Gravidade calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new row's cell values. The order below matters: it reproduces
# the order in which the workbook's author introduced new shared-string
# entries (G, H, I, then B, then J — the remaining cells reuse strings
# that already exist elsewhere in the sheet).
$ws.Range("A33").Value = 32
$ws.Range("G33").Value = "Homem Horas de Exposição ao Risco pode não corresponder com o número total de homem horas trabalhada."
$ws.Range("H33").Value = 'Implementar "Fator de Ajuste" para o calculo das Taxas de Frequência e Gravidade.'
$ws.Range("I33").Value = "Pendente (Testar)"
$ws.Range("B33").Value = "Pirelli"
$ws.Range("J33").Value = "Criar Variável Homem Hora Exposição ao Risco."
$ws.Range("C33").Value = "Processamento"
$ws.Range("D33").Value = "Taxas de Gravidade e Regressões"
$ws.Range("E33").Value = "PNL"
$ws.Range("F33").Value = "Mudança"
$ws.Range("K33").Value = "Não"

# Match the row height used for this new entry (wrapped text, 45pt tall).
$ws.Rows.Item(33).RowHeight = 45

# Move the active selection down to the newly added row, as the author
# left it positioned there after entering the data.
$ws.Range("A33").Select() | Out-Null
